$d = $word.ActiveDocument

# --- operation 1: replace before-paragraphs 85..85 ---
$r = $d.Paragraphs.Item(85).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>The data for the project is obtained from different csv files which are linked directly to the project. As a result, the Data gets Updated when the csv file is updated in the website.
        </w:t>
      </w:r>
    </w:p>'
$r.InsertXML($xml)

# --- operation 2: insert before-paragraphs 84..83 ---
$r = $d.Paragraphs.Item(84).Range
$r.Collapse(1)  # wdCollapseStart
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:b/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>'
$r.InsertXML($xml)

# --- operation 3: replace before-paragraphs 70..70 ---
$r = $d.Paragraphs.Item(70).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="12045" w:dyaOrig="6468" w14:anchorId="60528DA3">
          <v:rect id="rectole0000000011" o:spid="_x0000_i1036" style="width:518.4pt;height:239.05pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId27" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000011" DrawAspect="Content" ObjectID="_1658764928" r:id="rId28"/>
        </w:object>
      </w:r>
    </w:p>'
$r.InsertXML($xml)

# --- operation 4: replace before-paragraphs 68..68 ---
$r = $d.Paragraphs.Item(68).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="360"/>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="8104" w:dyaOrig="963" w14:anchorId="437AFD41">
          <v:rect id="rectole0000000010" o:spid="_x0000_i1035" style="width:406.1pt;height:48.95pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId25" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000010" DrawAspect="Content" ObjectID="_1658764927" r:id="rId26"/>
        </w:object>
      </w:r>
    </w:p>'
$r.InsertXML($xml)

# --- operation 5: replace before-paragraphs 62..65 ---
$r = $d.Range($d.Paragraphs.Item(62).Range.Start, $d.Paragraphs.Item(65).Range.End)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="720"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="8366" w:dyaOrig="907" w14:anchorId="12832DEE">
          <v:rect id="rectole0000000006" o:spid="_x0000_i1031" style="width:417.6pt;height:46.1pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId17" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000006" DrawAspect="Content" ObjectID="_1658764923" r:id="rId18"/>
        </w:object>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="12350" w:dyaOrig="6004" w14:anchorId="4D85B0A4">
          <v:rect id="rectole0000000007" o:spid="_x0000_i1032" style="width:506.9pt;height:3in" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId19" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000007" DrawAspect="Content" ObjectID="_1658764924" r:id="rId20"/>
        </w:object>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="360"/>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="8564" w:dyaOrig="964" w14:anchorId="1CD671BA">
          <v:rect id="rectole0000000008" o:spid="_x0000_i1033" style="width:429.1pt;height:48.95pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId21" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000008" DrawAspect="Content" ObjectID="_1658764925" r:id="rId22"/>
        </w:object>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:object w:dxaOrig="13161" w:dyaOrig="6370" w14:anchorId="69D5A07E">
          <v:rect id="rectole0000000009" o:spid="_x0000_i1034" style="width:509.75pt;height:233.3pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId23" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000009" DrawAspect="Content" ObjectID="_1658764926" r:id="rId24"/>
        </w:object>
      </w:r>
    </w:p>'
$r.InsertXML($xml)

# --- operation 6: replace before-paragraphs 47..48 ---
$r = $d.Range($d.Paragraphs.Item(47).Range.Start, $d.Paragraphs.Item(48).Range.End)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="360"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="4132" w:dyaOrig="993" w14:anchorId="64F056EF">
          <v:rect id="rectole0000000004" o:spid="_x0000_i1029" style="width:207.35pt;height:48.95pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId13" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000004" DrawAspect="Content" ObjectID="_1658764921" r:id="rId14"/>
        </w:object>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="10820" w:dyaOrig="6048" w14:anchorId="7343FE64">
          <v:rect id="rectole0000000005" o:spid="_x0000_i1030" style="width:501.1pt;height:273.6pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId15" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000005" DrawAspect="Content" ObjectID="_1658764922" r:id="rId16"/>
        </w:object>
      </w:r>
    </w:p>'
$r.InsertXML($xml)

# --- operation 7: insert before-paragraphs 44..43 ---
$r = $d.Paragraphs.Item(44).Range
$r.Collapse(1)  # wdCollapseStart
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>'
$r.InsertXML($xml)

# --- operation 8: replace before-paragraphs 30..30 ---
$r = $d.Paragraphs.Item(30).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="12282" w:dyaOrig="5904" w14:anchorId="222AF471">
          <v:rect id="rectole0000000003" o:spid="_x0000_i1028" style="width:509.75pt;height:204.5pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId11" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000003" DrawAspect="Content" ObjectID="_1658764920" r:id="rId12"/>
        </w:object>
      </w:r>
    </w:p>'
$r.InsertXML($xml)

# --- operation 9: insert before-paragraphs 27..26 ---
$r = $d.Paragraphs.Item(27).Range
$r.Collapse(1)  # wdCollapseStart
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>'
$r.InsertXML($xml)

# --- operation 10: replace before-paragraphs 25..25 ---
$r = $d.Paragraphs.Item(25).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="12361" w:dyaOrig="4982" w14:anchorId="55FB3B9E">
          <v:rect id="rectole0000000002" o:spid="_x0000_i1027" style="width:524.15pt;height:172.8pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId9" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000002" DrawAspect="Content" ObjectID="_1658764919" r:id="rId10"/>
        </w:object>
      </w:r>
    </w:p>'
$r.InsertXML($xml)

# --- operation 11: replace before-paragraphs 20..20 ---
$r = $d.Paragraphs.Item(20).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:ind w:left="-709"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:b/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="11664" w:dyaOrig="5005" w14:anchorId="52605DCC">
          <v:rect id="rectole0000000001" o:spid="_x0000_i1026" style="width:544.3pt;height:218.9pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId7" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000001" DrawAspect="Content" ObjectID="_1658764918" r:id="rId8"/>
        </w:object>
      </w:r>
    </w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>'
$r.InsertXML($xml)

# --- operation 12: replace before-paragraphs 14..14 ---
$r = $d.Paragraphs.Item(14).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:b/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>'
$r.InsertXML($xml)

# --- operation 13: replace before-paragraphs 3..3 ---
$r = $d.Paragraphs.Item(3).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:object w:dxaOrig="8908" w:dyaOrig="4998" w14:anchorId="279D234F">
          <v:rect id="rectole0000000000" o:spid="_x0000_i1025" style="width:446.4pt;height:250.55pt" o:ole="" o:preferrelative="t" stroked="f">
            <v:imagedata r:id="rId5" o:title=""/>
          </v:rect>
          <o:OLEObject Type="Embed" ProgID="StaticMetafile" ShapeID="rectole0000000000" DrawAspect="Content" ObjectID="_1658764917" r:id="rId6"/>
        </w:object>
      </w:r>
    </w:p>'
$r.InsertXML($xml)

Write-Output "done"